$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.348.72'
$ws.Range('D2').Style = $style
$ws.Range('E2').Value = '  -1.88%  '
$style = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.852.01'
$ws.Range('D3').Style = $style
$ws.Range('E3').Value = '  -1.27%  '
$style = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = $style
$ws.Range('E4').Value = '  +0.16%  '
$style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.6985'
$ws.Range('D5').Style = $style
$ws.Range('E5').Value = '  -5.75%  '
$style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '239.00'
$ws.Range('D6').Style = $style
$ws.Range('E6').Value = '  -1.49%  '
$ws.Range('E7').Value = '  +0.21%  '
$style = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3067'
$ws.Range('D8').Style = $style
$ws.Range('E8').Value = '  -2.83%  '
$style = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07386'
$ws.Range('D9').Style = $style
$ws.Range('E9').Value = '  +2.50%  '
$style = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.68'
$ws.Range('D10').Style = $style
$ws.Range('E10').Value = '  -4.14%  '
$style = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08110'
$ws.Range('D11').Style = $style
$ws.Range('E11').Value = '  -3.25%  '
$style = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.887.77'
$ws.Range('D12').Style = $style
$ws.Range('E12').Value = '  +1.27%  '
$style = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7247'
$ws.Range('D13').Style = $style
$ws.Range('E13').Value = '  -3.56%  '
$style = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.211'
$ws.Range('D14').Style = $style
$ws.Range('E14').Value = '  -4.11%  '
$style = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.56'
$ws.Range('D15').Style = $style
$ws.Range('E15').Value = '  -3.36%  '
$style = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.500.31'
$ws.Range('D16').Style = $style
$ws.Range('E16').Value = '  -1.39%  '
$style = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.901'
$ws.Range('D17').Style = $style
$ws.Range('E17').Value = '  -3.16%  '
$style = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '242.30'
$ws.Range('D18').Style = $style
$ws.Range('E18').Value = '  -1.79%  '
$style = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007717'
$ws.Range('D19').Style = $style
$ws.Range('E19').Value = '  -1.73%  '
$style = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.10'
$ws.Range('D20').Style = $style
$ws.Range('E20').Value = '  -3.72%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$style = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('D21').Style = $style
$ws.Range('E21').Value = '  +0.35%  '
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$style = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.135.26'
$ws.Range('D22').Style = $style
$ws.Range('E22').Value = '  +0.39%  '
$style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.002'
$ws.Range('D23').Style = $style
$ws.Range('E23').Value = '  +0.15%  '
$style = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.615'
$ws.Range('D24').Style = $style
$ws.Range('E24').Value = '  -5.04%  '
$style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1482'
$ws.Range('D25').Style = $style
$ws.Range('E25').Value = '  -4.61%  '
$style = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.027'
$ws.Range('D26').Style = $style
$ws.Range('E26').Value = '  -2.63%  '
$style = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '161.13'
$ws.Range('D27').Style = $style
$ws.Range('E27').Value = '  -2.45%  '
$style = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.06'
$ws.Range('D28').Style = $style
$ws.Range('E28').Value = '  -3.19%  '
$style = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.943'
$ws.Range('D29').Style = $style
$ws.Range('E29').Value = '  -4.72%  '
$ws.Range('E30').Value = '  -7.59%  '
$ws.Range('E31').Value = '  -1.74%  '
$style = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.409'
$ws.Range('D32').Style = $style
$ws.Range('E32').Value = '  -4.34%  '
$style = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.058'
$ws.Range('D33').Style = $style
$style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05270'
$ws.Range('D34').Style = $style
$ws.Range('E34').Value = '  -0.87%  '
$style = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.198'
$ws.Range('D35').Style = $style
$ws.Range('E35').Value = '  -3.28%  '
$style = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7223'
$ws.Range('D36').Style = $style
$ws.Range('E36').Value = '  -4.30%  '
$style = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.004'
$ws.Range('D37').Style = $style
$ws.Range('E37').Value = '  +0.41%  '
$style = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.682'
$ws.Range('D38').Style = $style
$ws.Range('E38').Value = '  -0.32%  '
$style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01865'
$ws.Range('D39').Style = $style
$ws.Range('E39').Value = '  -4.95%  '
$style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.724'
$ws.Range('D40').Style = $style
$ws.Range('E40').Value = '  -1.20%  '
$style = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8748'
$ws.Range('D41').Style = $style
$ws.Range('E41').Value = '  +2.18%  '
$style = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4309'
$ws.Range('D42').Style = $style
$ws.Range('E42').Value = '  -4.39%  '
$style = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.930'
$ws.Range('D43').Style = $style
$ws.Range('E43').Value = '  -2.24%  '
$style = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '69.79'
$ws.Range('D44').Style = $style
$ws.Range('E44').Value = '  -3.91%  '
$ws.Range('E45').Value = '  +0.06%  '
$style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.036.62'
$ws.Range('D46').Style = $style
$ws.Range('E46').Value = '  -6.53%  '
$style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.34'
$ws.Range('D47').Style = $style
$ws.Range('E47').Value = '  -0.79%  '
$style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.258'
$ws.Range('D48').Style = $style
$ws.Range('E48').Value = '  -4.89%  '
$style = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.030.71'
$ws.Range('D49').Style = $style
$ws.Range('E50').Value = '  -5.18%  '
$style = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.182'
$ws.Range('D51').Style = $style
$ws.Range('E51').Value = '  -3.21%  '
